$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.447.65"
$ws.Range("E2").Value = "  +1.34%  "

# Row 3
$ws.Range("D3").Value = "1.675.76"
$ws.Range("E3").Value = "  +2.07%  "

# Row 5
$ws.Range("D5").Value = "'219.39"
$ws.Range("E5").Value = "  +2.18%  "

# Row 6
$ws.Range("D6").Value = "'0.5311"
$ws.Range("E6").Value = "  +1.76%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.2697"
$ws.Range("E8").Value = "  +3.69%  "

# Row 9
$ws.Range("D9").Value = "'0.06395"
$ws.Range("E9").Value = "  +1.34%  "

# Row 10
$ws.Range("D10").Value = "'21.80"
$ws.Range("E10").Value = "  +5.36%  "

# Row 11
$ws.Range("D11").Value = "'0.07803"
$ws.Range("E11").Value = "  +1.58%  "

# Row 12
$ws.Range("D12").Value = "1.677.89"
$ws.Range("E12").Value = "  +2.29%  "

# Row 13
$ws.Range("D13").Value = "'4.509"
$ws.Range("E13").Value = "  +2.17%  "

# Row 14
$ws.Range("D14").Value = "'0.5583"
$ws.Range("E14").Value = "  +0.52%  "

# Row 15
$ws.Range("D15").Value = "0.0₅8333"
$ws.Range("E15").Value = "  +1.79%  "

# Row 16
$ws.Range("D16").Value = "'65.63"
$ws.Range("E16").Value = "  +0.87%  "

# Row 17
$ws.Range("D17").Value = "26.482.09"
$ws.Range("E17").Value = "  +1.48%  "

# Row 18
$ws.Range("E18").Value = "  -0.04%  "

# Row 19
$ws.Range("D19").Value = "'4.779"
$ws.Range("E19").Value = "  +1.35%  "

# Row 20
$ws.Range("D20").Value = "'192.90"
$ws.Range("E20").Value = "  +2.04%  "

# Row 21
$ws.Range("D21").Value = "'10.29"

# Row 22
$ws.Range("D22").Value = "'6.321"
$ws.Range("E22").Value = "  +2.37%  "

# Row 23
$ws.Range("E23").Value = "  +0.09%  "

# Row 24
$ws.Range("D24").Value = "'0.1276"
$ws.Range("E24").Value = "  +5.93%  "

# Row 25
$ws.Range("D25").Value = "'140.02"
$ws.Range("E25").Value = "  -3.73%  "

# Row 26
$ws.Range("D26").Value = "'7.404"
$ws.Range("E26").Value = "  -0.11%  "

# Row 27
$ws.Range("E27").Value = "  +2.84%  "

# Row 28
$ws.Range("D28").Value = "'1.444"
$ws.Range("E28").Value = "  +3.83%  "

# Row 29
$ws.Range("D29").Value = "'0.06263"
$ws.Range("E29").Value = "  +6.19%  "

# Row 30
$ws.Range("E30").Value = "  +2.09%  "

# Row 31
$ws.Range("D31").Value = "'3.608"
$ws.Range("E31").Value = "  +4.84%  "

# Row 32
$ws.Range("D32").Value = "'3.454"
$ws.Range("E32").Value = "  +1.44%  "

# Row 33
$ws.Range("D33").Value = "'1.692"
$ws.Range("E33").Value = "  +2.59%  "

# Row 34
$ws.Range("E34").Value = "  +2.77%  "

# Row 35
$ws.Range("D35").Value = "'0.6176"
$ws.Range("E35").Value = "  +9.24%  "

# Row 36
$ws.Range("D36").Value = "'2.422"
$ws.Range("E36").Value = "  +1.31%  "

# Row 37
$ws.Range("D37").Value = "'2.787"
$ws.Range("E37").Value = "  +1.17%  "

# Row 38
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "'6.165"
$ws.Range("E38").Value = "  +8.15%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01630"
$ws.Range("E39").Value = "  +0.81%  "

# Row 40
$ws.Range("D40").Value = "1.095.05"
$ws.Range("E40").Value = "  +6.62%  "

# Row 41
$ws.Range("D41").Value = "'0.8630"
$ws.Range("E41").Value = "  +1.14%  "

# Row 42
$ws.Range("E42").Value = "  -0.09%  "

# Row 43
$ws.Range("D43").Value = "'100.53"
$ws.Range("E43").Value = "  +0.38%  "

# Row 44
$ws.Range("E44").Value = "  +1.67%  "

# Row 45
$ws.Range("D45").Value = "0.0₈111"
$ws.Range("E45").Value = "  +3.61%  "

# Row 46
$ws.Range("D46").Value = "'58.71"
$ws.Range("E46").Value = "  +5.29%  "

# Row 47
$ws.Range("D47").Value = "'8.143"
$ws.Range("E47").Value = "  +0.66%  "

# Row 48
$ws.Range("D48").Value = "'1.004"
$ws.Range("E48").Value = "  +0.13%  "

# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.484"
$ws.Range("E49").Value = "  +7.09%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05199"
$ws.Range("E50").Value = "  +0.96%  "
